$wb = $excel.ActiveWorkbook

$wsEst = $wb.Worksheets.Item("Estimation Tool")
$wsVar = $wb.Worksheets.Item("Variables")

# --- Variables sheet: relabel rows (text content changes / reshuffle) ---
$wsVar.Range("E2").Value = "Power Apps per app pay-as-you-go"
$wsVar.Range("B3").Value = "Cost per user/month"
$wsVar.Range("B5").Value = "Default Dataverse Database capacity per tenant"
$wsVar.Range("B6").Value = "Default Dataverse File capacity per tenant"
$wsVar.Range("B7").Value = "Default Dataverse Log capacity per tenant"
$wsVar.Range("B9").Value = "Dataverse Database capacity accrued per user/month (GBs)"
$wsVar.Range("B10").Value = "Dataverse File capacity accrued per user/month (GBs)"
$wsVar.Range("B11").Value = "Dataverse Log capacity accrued per user/month (GBs)"
$wsVar.Range("B13").Value = "Dataverse Database Capacity add-on cost per month per GB"
$wsVar.Range("B14").Value = "Dataverse File Capacity add-on cost per month per GB"
$wsVar.Range("B15").Value = "Dataverse Log Capacity add-on cost per month per GB"
$wsVar.Range("B17").Value = "Allocated Requests per 24 Hours"

# --- Estimation Tool sheet: clear usage inputs ---
$wsEst.Range("C3").ClearContents()
$wsEst.Range("C4").ClearContents()
$wsEst.Range("C5").ClearContents()

# --- Header row 7: add labels ---
$wsEst.Range("B7").Value = "License Model"
$wsEst.Range("C7").Value = "Cost"
$wsEst.Rows.Item(7).RowHeight = 29.4

# --- Update selections to match final state ---
$wsVar.Range("B16").Select()
$wsEst.Range("G7").Select()

$wb.Save()
